$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 192.57143
$ws.Range("I12").Value = 89.59999999999999
$ws.Range("J12").Value = 450
$ws.Range("K12").Value = 89.59999999999999
$ws.Range("L12").Value = 450
$ws.Range("M12").Value = 80.40000000000001
$ws.Range("N12").Value = -790
$ws.Range("H15").Value = 350010.1
$ws.Range("I15").Value = 350010.1
$ws.Range("K15").Value = 1050030.3
$ws.Range("M15").Value = -1049861.3
$ws.Range("H40").Value = 4242
$ws.Range("I40").Value = 2000
$ws.Range("J40").Value = 4690.4
$ws.Range("K40").Value = 2000
$ws.Range("L40").Value = 4690.4
$ws.Range("M40").Value = -1825
$ws.Range("N40").Value = -5040.4
$ws.Range("H80").Value = 591.7368
$ws.Range("I80").Value = 381.625
$ws.Range("J80").Value = 744.5454999999999
$ws.Range("K80").Value = 1144.875
$ws.Range("L80").Value = 2233.6365
$ws.Range("M80").Value = -146.875
$ws.Range("N80").Value = -4229.6365
$ws.Range("H83").Value = 591.7368
$ws.Range("I83").Value = 381.625
$ws.Range("J83").Value = 744.5454999999999
$ws.Range("K83").Value = 3434.625
$ws.Range("L83").Value = 6700.9095
$ws.Range("M83").Value = 1557.375
$ws.Range("N83").Value = -16684.9095
$ws.Range("H88").Value = 1774.5
$ws.Range("J88").Value = 2062.9
$ws.Range("L88").Value = 2062.9
$ws.Range("N88").Value = -2874.9
$ws.Range("H91").Value = 1774.5
$ws.Range("J91").Value = 2062.9
$ws.Range("L91").Value = 2062.9
$ws.Range("N91").Value = -4870.9
$ws.Range("H106").Value = 2507.2144
$ws.Range("I106").Value = 828.0909
$ws.Range("K106").Value = 828.0909
$ws.Range("M106").Value = -197.0909
$ws.Range("H115").Value = 2741.5
$ws.Range("I115").Value = 1487
$ws.Range("J115").Value = 4832.3335
$ws.Range("K115").Value = 4461
$ws.Range("L115").Value = 14497.0005
$ws.Range("M115").Value = -2894
$ws.Range("N115").Value = -17631.0005
$ws.Range("H138").Value = 3042.9355
$ws.Range("I138").Value = 865.4138
$ws.Range("J138").Value = 4956.515
$ws.Range("K138").Value = 2596.2414
$ws.Range("L138").Value = 14869.545
$ws.Range("M138").Value = 2543.7586
$ws.Range("N138").Value = -25149.545

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3535.1482
$ws.Range("I2").Value = 2206.25
$ws.Range("K2").Value = 2206.25
$ws.Range("M2").Value = -2093.25
$ws.Range("H32").Value = 1250259.4
$ws.Range("I32").Value = 754
$ws.Range("K32").Value = 754
$ws.Range("M32").Value = -467
$ws.Range("H74").Value = 9186.75
$ws.Range("I74").Value = 11198.8
$ws.Range("K74").Value = 11198.8
$ws.Range("M74").Value = -10324.8
$ws.Range("H77").Value = 9186.75
$ws.Range("I77").Value = 11198.8
$ws.Range("K77").Value = 55994
$ws.Range("M77").Value = -51626
$ws.Range("H116").Value = 3535.1482
$ws.Range("I116").Value = 2206.25
$ws.Range("K116").Value = 2206.25
$ws.Range("M116").Value = 87.75
$ws.Range("H122").Value = 5600.3
$ws.Range("I122").Value = 4367.1665
$ws.Range("K122").Value = 13101.4995
$ws.Range("M122").Value = -10651.4995
$ws.Range("H132").Value = 1366085.9
$ws.Range("I132").Value = 1895957.9
$ws.Range("K132").Value = 5687873.699999999
$ws.Range("M132").Value = -5685343.699999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3535.1482
$ws.Range("I3").Value = 2206.25
$ws.Range("K3").Value = 2206.25
$ws.Range("M3").Value = -2092.25
$ws.Range("H74").Value = 27633.908
$ws.Range("J74").Value = 27633.908
$ws.Range("L74").Value = 27633.908
$ws.Range("N74").Value = -29505.908
$ws.Range("H77").Value = 27633.908
$ws.Range("J77").Value = 27633.908
$ws.Range("L77").Value = 82901.724
$ws.Range("N77").Value = -92261.724
$ws.Range("H93").Value = 55000
$ws.Range("J93").Value = 55000
$ws.Range("L93").Value = 55000
$ws.Range("N93").Value = -58744
$ws.Range("H97").Value = 20175.076
$ws.Range("I97").Value = 13534.625
$ws.Range("J97").Value = 30799.8
$ws.Range("K97").Value = 13534.625
$ws.Range("L97").Value = 30799.8
$ws.Range("M97").Value = -12543.625
$ws.Range("N97").Value = -32781.8
$ws.Range("H99").Value = 8369.351000000001
$ws.Range("I99").Value = 8190.543
$ws.Range("J99").Value = 8653.817999999999
$ws.Range("K99").Value = 8190.543
$ws.Range("L99").Value = 8653.817999999999
$ws.Range("M99").Value = -6692.543
$ws.Range("N99").Value = -11649.818

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 13891235
$ws.Range("I99").Value = 37038340
$ws.Range("J99").Value = 2972
$ws.Range("K99").Value = 37038340
$ws.Range("L99").Value = 2972
$ws.Range("M99").Value = -37036842
$ws.Range("N99").Value = -5968
$ws.Range("H122").Value = 3365.25
$ws.Range("I122").Value = 2635.5
$ws.Range("K122").Value = 7906.5
$ws.Range("M122").Value = -5456.5
$ws.Range("H126").Value = 13891235
$ws.Range("I126").Value = 37038340
$ws.Range("J126").Value = 2972
$ws.Range("K126").Value = 111115020
$ws.Range("L126").Value = 8916
$ws.Range("M126").Value = -111112550
$ws.Range("N126").Value = -13856

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 385709.84
$ws.Range("I5").Value = 753
$ws.Range("K5").Value = 2259
$ws.Range("M5").Value = -2147
$ws.Range("H37").Value = 135994.5
$ws.Range("J37").Value = 135994.5
$ws.Range("L37").Value = 407983.5
$ws.Range("N37").Value = -408207.5
$ws.Range("H94").Value = 5015.846
$ws.Range("J94").Value = 4642.9
$ws.Range("L94").Value = 13928.7
$ws.Range("N94").Value = -15280.7
$ws.Range("H116").Value = 1165.762
$ws.Range("I116").Value = 1111.1333
$ws.Range("K116").Value = 3333.3999
$ws.Range("M116").Value = 108.6001000000001
$ws.Range("H122").Value = 81791.14
$ws.Range("I122").Value = 499.53333
$ws.Range("J122").Value = 107194.77
$ws.Range("K122").Value = 4495.79997
$ws.Range("L122").Value = 964752.9300000001
$ws.Range("M122").Value = -2045.79997
$ws.Range("N122").Value = -969652.9300000001
$ws.Range("H125").Value = 13332
$ws.Range("J125").Value = 14998
$ws.Range("L125").Value = 44994
$ws.Range("N125").Value = -54834
$ws.Range("H130").Value = 19324.75
$ws.Range("I130").Value = 10000
$ws.Range("J130").Value = 22433
$ws.Range("K130").Value = 30000
$ws.Range("L130").Value = 67299
$ws.Range("M130").Value = -24980
$ws.Range("N130").Value = -77339
$ws.Range("H132").Value = 2400.611
$ws.Range("I132").Value = 1087.4445
$ws.Range("J132").Value = 3713.7778
$ws.Range("K132").Value = 9787.0005
$ws.Range("L132").Value = 33424.00019999999
$ws.Range("M132").Value = -7257.0005
$ws.Range("N132").Value = -38484.00019999999
$ws.Range("H134").Value = 14943.105
$ws.Range("I134").Value = 10227.934
$ws.Range("K134").Value = 30683.802
$ws.Range("M134").Value = -25613.802
$ws.Range("H135").Value = 385709.84
$ws.Range("I135").Value = 753
$ws.Range("K135").Value = 6777
$ws.Range("M135").Value = -4242
$ws.Range("H139").Value = 26318044
$ws.Range("I139").Value = 33335284
$ws.Range("J139").Value = 3392.5
$ws.Range("K139").Value = 100005852
$ws.Range("L139").Value = 10177.5
$ws.Range("M139").Value = -100000712
$ws.Range("N139").Value = -20457.5
$ws.Range("H140").Value = 34948190
$ws.Range("I140").Value = 47103010
$ws.Range("J140").Value = 3085.125
$ws.Range("K140").Value = 141309030
$ws.Range("L140").Value = 9255.375
$ws.Range("M140").Value = -141303850
$ws.Range("N140").Value = -19615.375

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 14590.36
$ws.Range("I70").Value = 9173.916999999999
$ws.Range("J70").Value = 19590.154
$ws.Range("K70").Value = 9173.916999999999
$ws.Range("L70").Value = 19590.154
$ws.Range("M70").Value = -8903.916999999999
$ws.Range("N70").Value = -20130.154
$ws.Range("H73").Value = 14590.36
$ws.Range("I73").Value = 9173.916999999999
$ws.Range("J73").Value = 19590.154
$ws.Range("K73").Value = 9173.916999999999
$ws.Range("L73").Value = 19590.154
$ws.Range("M73").Value = -8237.916999999999
$ws.Range("N73").Value = -21462.154
$ws.Range("H132").Value = 41668730
$ws.Range("I132").Value = 47620776
$ws.Range("K132").Value = 142862328
$ws.Range("M132").Value = -142859798

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 5879.6294
$ws.Range("I40").Value = 4736.5293
$ws.Range("J40").Value = 7822.9
$ws.Range("K40").Value = 4736.5293
$ws.Range("L40").Value = 7822.9
$ws.Range("M40").Value = -4600.5293
$ws.Range("N40").Value = -8094.9
$ws.Range("H55").Value = 2823.2432
$ws.Range("I55").Value = 1241.0741
$ws.Range("J55").Value = 7095.1
$ws.Range("K55").Value = 1241.0741
$ws.Range("L55").Value = 7095.1
$ws.Range("M55").Value = -1068.0741
$ws.Range("N55").Value = -7441.1
$ws.Range("H68").Value = 2554.9
$ws.Range("J68").Value = 2956.8572
$ws.Range("L68").Value = 2956.8572
$ws.Range("N68").Value = -4454.8572
$ws.Range("H71").Value = 2554.9
$ws.Range("J71").Value = 2956.8572
$ws.Range("L71").Value = 14784.286
$ws.Range("N71").Value = -22272.286
$ws.Range("H132").Value = 4999.0835
$ws.Range("I132").Value = 4089.9092
$ws.Range("K132").Value = 12269.7276
$ws.Range("M132").Value = -9739.7276
$ws.Range("H136").Value = 125006650
$ws.Range("I136").Value = 50007628
$ws.Range("K136").Value = 150022884
$ws.Range("M136").Value = -150020334

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 5747935.5
$ws.Range("I113").Value = 8772922
$ws.Range("J113").Value = 461.2
$ws.Range("K113").Value = 26318766
$ws.Range("L113").Value = 1383.6
$ws.Range("M113").Value = -26316596
$ws.Range("N113").Value = -5723.6
$ws.Range("H132").Value = 7750.88
$ws.Range("I132").Value = 4673.1577
$ws.Range("J132").Value = 17497
$ws.Range("K132").Value = 14019.4731
$ws.Range("L132").Value = 52491
$ws.Range("M132").Value = -11489.4731
$ws.Range("N132").Value = -57551

Write-Host "Updated 259 cells across 8 sheets"